# Auto-generated edit script applying cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.848.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "'3.439.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'583.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "'173.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'3.438.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D11").Value = "'6.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").Value = "'0.410"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").Value = "'4.039.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "'28.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.65%  "
$ws.Range("D16").Value = "'65.883.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "'0.0000170"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "'3.452.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "'5.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Value = "'13.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "'369.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").Value = "'7.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").Value = "'72.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "'0.528"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "'0.0000121"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.19%  "
$ws.Range("D27").Value = "'9.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("E28").Value = "  +3.72%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").Value = "'23.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("D37").Value = "'160.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "'0.878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").Value = "'28.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.01%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").Value = "'2.761.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").Value = "'4.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D46").Value = "'40.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "'24.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "'323.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "'6.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.83%  "
